$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the BOM date (row 3) ---
$ws.Range("A3").Value = "1/15/2013"

# --- Remove the old U3 / voltage-reg row (row 22) and insert the new
#     USB-receptacle row (K1, K3) in its place at row 21 ---
$ws.Rows("22:22").Delete()
$ws.Rows("21:21").Insert()

$ws.Range("C21").Value = "K1, K3"
$ws.Range("D21").Value = "USB receptacle"
$ws.Range("E21").Value = "NA"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = "USB receptacle"
$ws.Range("G21").Value = "CNC Tech"
$ws.Range("H21").Value = "1002-001-01000"
$ws.Range("I21").Value = "Digi-Key"
$ws.Range("J21").Value = "1175-1015-ND"
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.69
$ws.Range("M21").Formula = "=K21*L21"

$ws.Hyperlinks.Add($ws.Range("J21"), "http://www.digikey.com/product-detail/en/1002-001-01000/1175-1015-ND/2396816")
$ws.Hyperlinks.Add($ws.Range("G21"), "http://digikey.com/Suppliers/us/CNC-Tech.page?lang=en")

# --- Footnote explaining "(alt)" parts ---
$ws.Range("C25").Value = "(alt) means it has the same footprint and we can test it as an alternative part"

# --- Widen the Vendor Part # column and tidy up the view ---
$ws.Columns("H:H").ColumnWidth = 23.67
$ws.Range("C3").Select()
